$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new tracker rows (68, 69) following the existing pattern.
$ws.Range("A68").Value = "G1"
$ws.Range("B68").Value = "Test1"
$ws.Range("C68").Value = 45894
$ws.Range("C68").NumberFormat = "YYYY-MM-DD"
$ws.Range("D68").Value = 0.7273041052711734
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = -0.01

$ws.Range("A69").Value = "G2"
$ws.Range("B69").Value = "sedrftgyhuioygtfrd"
$ws.Range("C69").Value = 45894
$ws.Range("C69").NumberFormat = "YYYY-MM-DD"
$ws.Range("D69").Value = 0.7273041052711734
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = -0.01
